$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, copy the style of row 16 (A16, style index 1) down to the three new rows (17-19) ---
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

# --- Rows 10-16: re-label & re-populate with the new (shifted) data ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.781294520948619
$ws.Range("D10").Value = 1.308354768009678
$ws.Range("E10").Value = 0.9245862354206003
$ws.Range("F10").Value = 1.781294520948619
$ws.Range("G10").Value = 0.9667450792716777
$ws.Range("H10").Value = 1.036695227981307
$ws.Range("I10").Value = 0.9440024818372539
$ws.Range("J10").Value = 1.308354768009678
$ws.Range("K10").Value = 1.116470501715139
$ws.Range("L10").Value = 1.448882511331879
$ws.Range("M10").Value = 1.160279718911523

$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.8773757491072123
$ws.Range("D11").Value = 1.278215570942997
$ws.Range("E11").Value = 1.203882458872939
$ws.Range("F11").Value = 0.8773757491072123
$ws.Range("G11").Value = 0.8259207892945777
$ws.Range("H11").Value = 1.92714238098915
$ws.Range("I11").Value = 0.8578269092266989
$ws.Range("J11").Value = 1.278215570942997
$ws.Range("K11").Value = 1.241049014907968
$ws.Range("L11").Value = 1.05921238200759
$ws.Range("M11").Value = 1.161727309738929

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.8742189613674587
$ws.Range("D12").Value = 1.280108494192596
$ws.Range("E12").Value = 1.204928734832489
$ws.Range("F12").Value = 0.8742189613674587
$ws.Range("G12").Value = 0.826826436576325
$ws.Range("H12").Value = 1.927430276158799
$ws.Range("I12").Value = 0.8571935061697218
$ws.Range("J12").Value = 1.280108494192596
$ws.Range("K12").Value = 1.242518614512542
$ws.Range("L12").Value = 1.058368787940001
$ws.Range("M12").Value = 1.161784401549565

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.8766649806272484
$ws.Range("D13").Value = 1.278302866677423
$ws.Range("E13").Value = 1.204538471416323
$ws.Range("F13").Value = 0.8766649806272484
$ws.Range("G13").Value = 0.8260818684702093
$ws.Range("H13").Value = 1.926757393121276
$ws.Range("I13").Value = 0.8574713656875421
$ws.Range("J13").Value = 1.278302866677423
$ws.Range("K13").Value = 1.241420669046873
$ws.Range("L13").Value = 1.059042824837061
$ws.Range("M13").Value = 1.16163615766667

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.3936919999999989
$ws.Range("D14").Value = 2.092115999999994
$ws.Range("E14").Value = 1.514191999999997
$ws.Range("F14").Value = 0.3936919999999989
$ws.Range("G14").Value = 0.7221799999999982
$ws.Range("H14").Value = 2.71490800000001
$ws.Range("I14").Value = 0.727136000000001
$ws.Range("J14").Value = 2.092115999999994
$ws.Range("K14").Value = 1.803153999999995
$ws.Range("L14").Value = 1.098422999999997
$ws.Range("M14").Value = 1.360704

$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.28
$ws.Range("D15").Value = 2.427662499999996
$ws.Range("E15").Value = 1.685850000000001
$ws.Range("F15").Value = 0.28
$ws.Range("G15").Value = 0.5423624999999999
$ws.Range("H15").Value = 3.286662500000002
$ws.Range("I15").Value = 0.6899999999999999
$ws.Range("J15").Value = 2.427662499999996
$ws.Range("K15").Value = 2.056756249999999
$ws.Range("L15").Value = 1.168378124999999
$ws.Range("M15").Value = 1.485422916666666

$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.5666026567680008
$ws.Range("D16").Value = 1.780130991308805
$ws.Range("E16").Value = 1.373312286720004
$ws.Range("F16").Value = 0.5666026567680008
$ws.Range("G16").Value = 0.7497090598911992
$ws.Range("H16").Value = 2.259105566617599
$ws.Range("I16").Value = 0.8239645196288005
$ws.Range("J16").Value = 1.780130991308805
$ws.Range("K16").Value = 1.576721639014405
$ws.Range("L16").Value = 1.071662147891203
$ws.Range("M16").Value = 1.258804180155735

# --- New rows 17-19 ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9861578971481871
$ws.Range("D17").Value = 0.9913306608720435
$ws.Range("E17").Value = 0.9881494912983582
$ws.Range("F17").Value = 0.9861578971481871
$ws.Range("G17").Value = 0.9890972860804736
$ws.Range("H17").Value = 0.9890787327105353
$ws.Range("I17").Value = 0.9863806546095567
$ws.Range("J17").Value = 0.9913306608720435
$ws.Range("K17").Value = 0.9897400760852009
$ws.Range("L17").Value = 0.9879489866166939
$ws.Range("M17").Value = 0.988365787119859

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9897540132686324
$ws.Range("D18").Value = 1.001310153597113
$ws.Range("E18").Value = 1.090295654112694
$ws.Range("F18").Value = 0.9897540132686324
$ws.Range("G18").Value = 0.9506232175549073
$ws.Range("H18").Value = 1.092845141924159
$ws.Range("I18").Value = 0.9468405064753135
$ws.Range("J18").Value = 1.001310153597113
$ws.Range("K18").Value = 1.045802903854903
$ws.Range("L18").Value = 1.017778458561768
$ws.Range("M18").Value = 1.01194478115547

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9942783707768669
$ws.Range("D19").Value = 0.9876431694675819
$ws.Range("E19").Value = 0.9801133740497162
$ws.Range("F19").Value = 0.9942783707768669
$ws.Range("G19").Value = 0.9881620024722707
$ws.Range("H19").Value = 0.9631200627381333
$ws.Range("I19").Value = 0.9859554117093244
$ws.Range("J19").Value = 0.9876431694675819
$ws.Range("K19").Value = 0.9838782717586491
$ws.Range("L19").Value = 0.989078321267758
$ws.Range("M19").Value = 0.9832120652023155

